# This workbook ("Hortaliza, Vega Central Mapocho de Santiago - Zapallo")
# gets a new weekly data block for the "Camote" variety, date 2021-11-22
# (serial 44522), inserted right before the existing 2021-02-19 block
# (which currently starts at row 650). Inserting the two new rows pushes
# all subsequent rows down by two, growing the sheet from A1:R747 to
# A1:R749 - matching the rest of the diff, where every later row's
# content equals the content that used to be two rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 650, shifting everything below down.
$ws.Rows.Item(650).Resize(2).Insert()

# New row 650: "1a nueva(o)" quality, origin Peru.
$ws.Cells.Item(650, 1).Value  = 9
$ws.Cells.Item(650, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(650, 3).Value  = "Metropolitana"
$ws.Cells.Item(650, 4).Value  = 44522
$ws.Cells.Item(650, 5).Value  = 13
$ws.Cells.Item(650, 6).Value  = 100112045
$ws.Cells.Item(650, 7).Value  = "Zapallo"
$ws.Cells.Item(650, 8).Value  = "Camote"
$ws.Cells.Item(650, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(650, 10).Value = 160
$ws.Cells.Item(650, 11).Value = 600
$ws.Cells.Item(650, 12).Value = 700
$ws.Cells.Item(650, 13).Value = 650
$ws.Cells.Item(650, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(650, 15).Value = "Perú"
$ws.Cells.Item(650, 16).Value = 650
$ws.Cells.Item(650, 17).Value = 1
$ws.Cells.Item(650, 18).Value = "Hortaliza"

# New row 651: "2a nueva(o)" quality, origin Peru.
$ws.Cells.Item(651, 1).Value  = 9
$ws.Cells.Item(651, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(651, 3).Value  = "Metropolitana"
$ws.Cells.Item(651, 4).Value  = 44522
$ws.Cells.Item(651, 5).Value  = 13
$ws.Cells.Item(651, 6).Value  = 100112045
$ws.Cells.Item(651, 7).Value  = "Zapallo"
$ws.Cells.Item(651, 8).Value  = "Camote"
$ws.Cells.Item(651, 9).Value  = "2a nueva(o)"
$ws.Cells.Item(651, 10).Value = 61
$ws.Cells.Item(651, 11).Value = 450
$ws.Cells.Item(651, 12).Value = 500
$ws.Cells.Item(651, 13).Value = 475
$ws.Cells.Item(651, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(651, 15).Value = "Perú"
$ws.Cells.Item(651, 16).Value = 475
$ws.Cells.Item(651, 17).Value = 1
$ws.Cells.Item(651, 18).Value = "Hortaliza"
